$d = $word.ActiveDocument

# 1) Remove the parenthetical "Domain CTO" nickname from the job title line.
$d.Content.Find.Execute(
    "Domain Technology Officer ("+[char]8220+"Domain CTO"+[char]8221+") Marketing for",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Domain Technology Officer Marketing for", 2)

# 2) Replace the quoted "CTO" nickname with "Technology Leader" in the body text.
$d.Content.Find.Execute(
    [char]8220+"CTO"+[char]8221+" of the Flixbus Marketing Organization",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Technology Leader of the Flixbus Marketing Organization", 2)
